$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 93 (pushes old rows 93-113 down to 94-114) ---
$ws.Rows.Item(93).Insert()

# Copy formatting (style) for the date cell in column A from the row below (now A94)
$ws.Range("A94").Copy()
$ws.Range("A93").PasteSpecial(-4122)

# --- 2. Fill in the values for the newly inserted row 93 ---
$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 2
$ws.Range("C93").Value = 8
$ws.Range("D93").Value = 174.4820065430752

# --- 3. Update the recalculated "somma mobile 7gg." (C) and "per 100mila abitanti" (D)
#         columns for the rows that were shifted down (now rows 90-112) ---
$ws.Range("C90").Value = 8
$ws.Range("D90").Value = 174.4820065430752

$ws.Range("C91").Value = 8
$ws.Range("D91").Value = 174.4820065430752

$ws.Range("C92").Value = 8
$ws.Range("D92").Value = 174.4820065430752

$ws.Range("C94").Value = 7
$ws.Range("D94").Value = 152.6717557251908

$ws.Range("C95").Value = 6
$ws.Range("D95").Value = 130.8615049073064

$ws.Range("C96").Value = 2
$ws.Range("D96").Value = 43.62050163576881

$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 0

$ws.Range("C98").Value = 1
$ws.Range("D98").Value = 21.81025081788441

$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 21.81025081788441

$ws.Range("C100").Value = 3
$ws.Range("D100").Value = 65.43075245365321

$ws.Range("C101").Value = 4
$ws.Range("D101").Value = 87.24100327153762

$ws.Range("C102").Value = 4
$ws.Range("D102").Value = 87.24100327153762

$ws.Range("C103").Value = 5
$ws.Range("D103").Value = 109.051254089422

$ws.Range("C104").Value = 6
$ws.Range("D104").Value = 130.8615049073064

$ws.Range("C105").Value = 10
$ws.Range("D105").Value = 218.1025081788441

$ws.Range("C106").Value = 10
$ws.Range("D106").Value = 218.1025081788441

$ws.Range("C107").Value = 8
$ws.Range("D107").Value = 174.4820065430752

$ws.Range("C108").Value = 8
$ws.Range("D108").Value = 174.4820065430752

$ws.Range("C109").Value = 9
$ws.Range("D109").Value = 196.2922573609597

$ws.Range("C110").Value = 9
$ws.Range("D110").Value = 196.2922573609597

$ws.Range("C111").Value = 12
$ws.Range("D111").Value = 261.7230098146129

$ws.Range("C112").Value = 7
$ws.Range("D112").Value = 152.6717557251908

# --- 4. Update the "nuovi pos." (B) values for the rows that changed relative counts ---
$ws.Range("B100").Value = 0
$ws.Range("B101").Value = 1
$ws.Range("B102").Value = 0
$ws.Range("B103").Value = 2
$ws.Range("B104").Value = 1
$ws.Range("B105").Value = 0
$ws.Range("B107").Value = 1
$ws.Range("B108").Value = 5
$ws.Range("B110").Value = 0

# --- 5. Append two brand-new rows (114, 115) at the end, copying the formatting
#         of row 113's date cell (column A) ---
$ws.Range("A113").Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)

$ws.Range("A114").Value = 44256
$ws.Range("B114").Value = 4
$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 0

# New rows have no prior C/D "blank" cells to inherit, so create them explicitly
# as empty-string cells to match the rest of the column's blank-value representation.
$ws.Range("C115").Value = ""
$ws.Range("D115").Value = ""

Write-Host "Final UsedRange:" $ws.UsedRange.Address()
